$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''317.99'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''3.66%'
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''39.94'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''2.69%'
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''5.143'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''0.80%'
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''0.08233'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''1.83%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''2.088'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''7.40%'
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''8.337'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''4.63%'
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9412'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''1.10%'
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.1370'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''-7.49%'
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1975'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''1.79%'
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.09111'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''-0.76%'
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03511'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''0.26%'
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09818'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''0.37%'
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001369'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''-2.14%'
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006341'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''7.28%'
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.696'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''-2.41%'
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''4.328'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''3.38%'
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''-4.96%'
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''0.3500'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''1.18%'
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''0.1311'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''0.69%'
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''4.971'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''9.24%'
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.2443'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''1.21%'
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''0.04343'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''-0.66%'
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''0.001223'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''-1.11%'
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.004842'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''13.11%'
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''-0.54%'
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''0.0003988'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''-10.34%'
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = '''0.02209'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''8.25%'
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.05216'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''2.65%'
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.007669'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''1.90%'
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.009653'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''-6.28%'
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''4.37%'
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.002032'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''-4.21%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.008883'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''-2.19%'
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.00006616'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''6.91%'
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.00000000748'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''-0.38%'
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''0.002926'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = '''0.001685'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''5.20%'
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.00002094'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''-0.38%'
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''0.0001994'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''-0.38%'
$ws.Range("E51").Style = "Normal"
